$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 33078.25
$ws.Range("J3").Value = 33078.25
$ws.Range("L3").Value = 33078.25
$ws.Range("N3").Value = -33306.25
$ws.Range("H17").Value = 878.0909
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 878.0909
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2634.2727
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = -2970.2727
$ws.Range("H102").Value = 33078.25
$ws.Range("J102").Value = 33078.25
$ws.Range("L102").Value = 33078.25
$ws.Range("N102").Value = -39568.25
$ws.Range("H103").Value = 487.63635
$ws.Range("I103").Value = 466
$ws.Range("J103").Value = 491.05264
$ws.Range("K103").Value = 1398
$ws.Range("L103").Value = 1473.15792
$ws.Range("M103").Value = -812
$ws.Range("N103").Value = -2645.15792
$ws.Range("H113").Value = 37446.43
$ws.Range("I113").Value = 92641.82000000001
$ws.Range("J113").Value = 1731.7646
$ws.Range("K113").Value = 92641.82000000001
$ws.Range("L113").Value = 1731.7646
$ws.Range("M113").Value = -89387.82000000001
$ws.Range("N113").Value = -8239.7646

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 253745
$ws.Range("I45").Value = 502000
$ws.Range("J45").Value = 5490
$ws.Range("K45").Value = 502000
$ws.Range("L45").Value = 5490
$ws.Range("M45").Value = -501623
$ws.Range("N45").Value = -6244
$ws.Range("H53").Value = 11510.75
$ws.Range("J53").Value = 11510.75
$ws.Range("L53").Value = 11510.75
$ws.Range("N53").Value = -12874.75
$ws.Range("H61").Value = 1809.0667
$ws.Range("I61").Value = 1653.0834
$ws.Range("K61").Value = 1653.0834
$ws.Range("M61").Value = -1441.0834
$ws.Range("H63").Value = 3515
$ws.Range("I63").Value = 2552.5
$ws.Range("J63").Value = 3900
$ws.Range("K63").Value = 2552.5
$ws.Range("L63").Value = 3900
$ws.Range("M63").Value = -1866.5
$ws.Range("N63").Value = -5272
$ws.Range("H66").Value = 3515
$ws.Range("I66").Value = 2552.5
$ws.Range("J66").Value = 3900
$ws.Range("K66").Value = 12762.5
$ws.Range("L66").Value = 19500
$ws.Range("M66").Value = -9330.5
$ws.Range("N66").Value = -26364
$ws.Range("H74").Value = 1842.2759
$ws.Range("J74").Value = 2896
$ws.Range("L74").Value = 2896
$ws.Range("N74").Value = -4644
$ws.Range("H77").Value = 1842.2759
$ws.Range("J77").Value = 2896
$ws.Range("L77").Value = 14480
$ws.Range("N77").Value = -23216
$ws.Range("H102").Value = 128377.375
$ws.Range("I102").Value = 335653.34
$ws.Range("J102").Value = 4011.8
$ws.Range("K102").Value = 335653.34
$ws.Range("L102").Value = 4011.8
$ws.Range("M102").Value = -334031.34
$ws.Range("N102").Value = -7255.8
$ws.Range("H132").Value = 1621.963
$ws.Range("I132").Value = 1241.4584
$ws.Range("K132").Value = 3724.3752
$ws.Range("M132").Value = -1194.3752
$ws.Range("H136").Value = 1809.0667
$ws.Range("I136").Value = 1653.0834
$ws.Range("K136").Value = 4959.2502
$ws.Range("M136").Value = -2409.2502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 857.4
$ws.Range("I94").Value = 715.6
$ws.Range("J94").Value = 999.2
$ws.Range("K94").Value = 715.6
$ws.Range("L94").Value = 999.2
$ws.Range("M94").Value = -264.6
$ws.Range("N94").Value = -1901.2
$ws.Range("H99").Value = 1498.4
$ws.Range("I99").Value = 1418.5555
$ws.Range("J99").Value = 1618.1666
$ws.Range("K99").Value = 1418.5555
$ws.Range("L99").Value = 1618.1666
$ws.Range("M99").Value = 79.44450000000006
$ws.Range("N99").Value = -4614.1666
$ws.Range("H134").Value = 2427.1875
$ws.Range("I134").Value = 2358.5334
$ws.Range("J134").Value = 3457
$ws.Range("K134").Value = 7075.600199999999
$ws.Range("L134").Value = 10371
$ws.Range("M134").Value = -4540.600199999999
$ws.Range("N134").Value = -15441

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 27209
$ws.Range("J120").Value = 27209
$ws.Range("L120").Value = 27209
$ws.Range("N120").Value = -34467

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 73133.64
$ws.Range("I119").Value = 100477.5
$ws.Range("J119").Value = 4774
$ws.Range("K119").Value = 301432.5
$ws.Range("L119").Value = 14322
$ws.Range("M119").Value = -296594.5
$ws.Range("N119").Value = -23998
$ws.Range("H132").Value = 2088.625
$ws.Range("I132").Value = 876
$ws.Range("J132").Value = 3301.25
$ws.Range("K132").Value = 7884
$ws.Range("L132").Value = 29711.25
$ws.Range("M132").Value = -5354
$ws.Range("N132").Value = -34771.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 166835310
$ws.Range("I80").Value = 250251970
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 250251970
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -250250972
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 166835310
$ws.Range("I83").Value = 250251970
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 1251259850
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -1251254858
$ws.Range("N83").Value = -19984
$ws.Range("H97").Value = 200004270
$ws.Range("I97").Value = 200004270
$ws.Range("K97").Value = 200004270
$ws.Range("M97").Value = -200003774
$ws.Range("H120").Value = 35462.6
$ws.Range("J120").Value = 35462.6
$ws.Range("L120").Value = 35462.6
$ws.Range("N120").Value = -45138.6
$ws.Range("H132").Value = 2910.2666
$ws.Range("I132").Value = 2534.5715
$ws.Range("J132").Value = 3239
$ws.Range("K132").Value = 7603.7145
$ws.Range("L132").Value = 9717
$ws.Range("M132").Value = -5073.7145
$ws.Range("N132").Value = -14777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 13999
$ws.Range("I11").Value = 13999
$ws.Range("K11").Value = 13999
$ws.Range("M11").Value = -13859
$ws.Range("H46").Value = 562868.4399999999
$ws.Range("I46").Value = 252.22223
$ws.Range("J46").Value = 1125484.6
$ws.Range("K46").Value = 252.22223
$ws.Range("L46").Value = 1125484.6
$ws.Range("M46").Value = -64.22223
$ws.Range("N46").Value = -1125860.6
$ws.Range("H100").Value = 1309.8667
$ws.Range("I100").Value = 1137.5
$ws.Range("J100").Value = 1372.5454
$ws.Range("K100").Value = 1137.5
$ws.Range("L100").Value = 1372.5454
$ws.Range("M100").Value = -596.5
$ws.Range("N100").Value = -2454.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 20063
$ws.Range("I6").Value = 310.25
$ws.Range("J6").Value = 46400
$ws.Range("K6").Value = 310.25
$ws.Range("L6").Value = 46400
$ws.Range("M6").Value = -195.25
$ws.Range("N6").Value = -46630
$ws.Range("H21").Value = 70017
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 70017
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 70017
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = -70487
$ws.Range("H29").Value = 70011
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = ""
$ws.Range("H35").Value = 70017
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 70017
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 70017
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = -70597
$ws.Range("H56").Value = 36485.273
$ws.Range("J56").Value = 39733.8
$ws.Range("L56").Value = 39733.8
$ws.Range("N56").Value = -41161.8
